# Update the "dSF" column (column F) values for a handful of rows to reflect
# the repulled/pushed data and recalculated mean.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    2  = -5
    5  = 3
    6  = -8
    13 = -2
    22 = -4
    24 = 2
    38 = -2
    45 = 0
    49 = 3
    55 = -4
    64 = -4
    65 = 4
}

foreach ($row in $updates.Keys) {
    $ws.Range("F$row").Value = $updates[$row]
}
